$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(101; 9; 30; 15; 60; 15)
    3 = @(902; 1; 0; 0; 0; 0)
    4 = @(901; 16; 15; 45; 60; 60)
    5 = @(1001; 18; 30; 75; 60; 72)
    6 = @(601; 9; 60; 67; 60; 42)
    7 = @(201; 9; 30; 15; 45; 30)
    8 = @(1203; 3; 15; 15; 15; 15)
    9 = @(501; 9; 52; 30; 75; 45)
    10 = @(701; 3; 90; 45; 97; 15)
    11 = @(801; 3; 67; 65; 52; 45)
    12 = @(301; 6; 45; 30; 60; 45)
    13 = @(401; 9; 48; 67; 75; 45)
    14 = @(1201; 2; 10; 10; 10; 10)
    15 = @(1202; 2; 10; 10; 10; 10)
    16 = @(502; 0; 4; 0; 0; 0)
    17 = @(802; 0; 4; 5; 4; 0)
    18 = @(3; 0; 3; 3; 3; 3)
    19 = @(1101; 0; 15; 30; 30; 0)
    20 = @(1; 0; 2; 2; 2; 2)
    21 = @(2; 0; 2; 2; 2; 2)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
